$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("AMSIN")

# Row 143 (existing tail row - refresh value/style)
$ws1.Cells.Item(143,1).ClearContents()
$ws1.Cells.Item(143,1).NumberFormat = "@"
$ws1.Cells.Item(143,1).Value = '2024-03-14'
$ws1.Cells.Item(143,1).Style = "Normal"
$ws1.Cells.Item(143,2).Value = 45365.55358833334
$ws1.Cells.Item(143,3).ClearContents()
$ws1.Cells.Item(143,3).Value = '189retest'
$ws1.Cells.Item(143,4).ClearContents()
$ws1.Cells.Item(143,4).Value = 269
$ws1.Cells.Item(143,5).ClearContents()
$ws1.Cells.Item(143,5).Value = 269
$ws1.Cells.Item(143,6).ClearContents()
$ws1.Cells.Item(143,6).Value = 0
$ws1.Cells.Item(143,7).ClearContents()
$ws1.Cells.Item(143,7).Value = 4.51

# Row 144
$ws1.Cells.Item(144,1).NumberFormat = "@"
$ws1.Cells.Item(144,1).Value = '2024-03-28'
$ws1.Cells.Item(144,1).Style = "Normal"
$ws1.Cells.Item(144,2).Style = "Normal"
$ws1.Cells.Item(144,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(144,2).Value = 45379.50678019676
$ws1.Cells.Item(144,3).Value = '190masstrail'
$ws1.Cells.Item(144,4).Value = 269
$ws1.Cells.Item(144,5).Value = 266
$ws1.Cells.Item(144,6).Value = 3
$ws1.Cells.Item(144,7).Value = 5.21

# Row 145
$ws1.Cells.Item(145,1).NumberFormat = "@"
$ws1.Cells.Item(145,1).Value = '2024-03-28'
$ws1.Cells.Item(145,1).Style = "Normal"
$ws1.Cells.Item(145,2).Style = "Normal"
$ws1.Cells.Item(145,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(145,2).Value = 45379.67118087963
$ws1.Cells.Item(145,3).Value = '190fstcycle'
$ws1.Cells.Item(145,4).Value = 269
$ws1.Cells.Item(145,5).Value = 268
$ws1.Cells.Item(145,6).Value = 1
$ws1.Cells.Item(145,7).Value = 5.21

# Row 146
$ws1.Cells.Item(146,1).NumberFormat = "@"
$ws1.Cells.Item(146,1).Value = '2024-03-29'
$ws1.Cells.Item(146,1).Style = "Normal"
$ws1.Cells.Item(146,2).Style = "Normal"
$ws1.Cells.Item(146,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(146,2).Value = 45380.45194586806
$ws1.Cells.Item(146,3).Value = '190scndcycle'
$ws1.Cells.Item(146,4).Value = 269
$ws1.Cells.Item(146,5).Value = 269
$ws1.Cells.Item(146,6).Value = 0
$ws1.Cells.Item(146,7).Value = 4.1

# Row 147
$ws1.Cells.Item(147,1).NumberFormat = "@"
$ws1.Cells.Item(147,1).Value = '2024-04-01'
$ws1.Cells.Item(147,1).Style = "Normal"
$ws1.Cells.Item(147,2).Style = "Normal"
$ws1.Cells.Item(147,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(147,2).Value = 45383.37318597222
$ws1.Cells.Item(147,3).Value = '190fnlrun'
$ws1.Cells.Item(147,4).Value = 269
$ws1.Cells.Item(147,5).Value = 269
$ws1.Cells.Item(147,6).Value = 0
$ws1.Cells.Item(147,7).Value = 4

# Row 148
$ws1.Cells.Item(148,1).NumberFormat = "@"
$ws1.Cells.Item(148,1).Value = '2024-04-25'
$ws1.Cells.Item(148,1).Style = "Normal"
$ws1.Cells.Item(148,2).Style = "Normal"
$ws1.Cells.Item(148,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(148,2).Value = 45407.45201841435
$ws1.Cells.Item(148,3).Value = '191trail'
$ws1.Cells.Item(148,4).Value = 269
$ws1.Cells.Item(148,5).Value = 266
$ws1.Cells.Item(148,6).Value = 3
$ws1.Cells.Item(148,7).Value = 5.18

# Row 149
$ws1.Cells.Item(149,1).NumberFormat = "@"
$ws1.Cells.Item(149,1).Value = '2024-05-02'
$ws1.Cells.Item(149,1).Style = "Normal"
$ws1.Cells.Item(149,2).Style = "Normal"
$ws1.Cells.Item(149,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(149,2).Value = 45414.48380873843
$ws1.Cells.Item(149,3).Value = '191fstcycle'
$ws1.Cells.Item(149,4).Value = 269
$ws1.Cells.Item(149,5).Value = 269
$ws1.Cells.Item(149,6).Value = 0
$ws1.Cells.Item(149,7).Value = 4.51

# Row 150
$ws1.Cells.Item(150,1).NumberFormat = "@"
$ws1.Cells.Item(150,1).Value = '2024-05-03'
$ws1.Cells.Item(150,1).Style = "Normal"
$ws1.Cells.Item(150,2).Style = "Normal"
$ws1.Cells.Item(150,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(150,2).Value = 45415.3477190625
$ws1.Cells.Item(150,3).Value = '191lstrun'
$ws1.Cells.Item(150,4).Value = 269
$ws1.Cells.Item(150,5).Value = 267
$ws1.Cells.Item(150,6).Value = 2
$ws1.Cells.Item(150,7).Value = 4.64

# Row 151
$ws1.Cells.Item(151,1).NumberFormat = "@"
$ws1.Cells.Item(151,1).Value = '2024-05-06'
$ws1.Cells.Item(151,1).Style = "Normal"
$ws1.Cells.Item(151,2).Style = "Normal"
$ws1.Cells.Item(151,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(151,2).Value = 45418.45908497685
$ws1.Cells.Item(151,3).Value = 'sampletestevent12'
$ws1.Cells.Item(151,4).Value = 269
$ws1.Cells.Item(151,5).Value = 269
$ws1.Cells.Item(151,6).Value = 0
$ws1.Cells.Item(151,7).Value = 4.13

# Row 152 (last row - default/unstyled like the old last row)
$ws1.Cells.Item(152,1).NumberFormat = "@"
$ws1.Cells.Item(152,1).Value = '2024-05-06'
$ws1.Cells.Item(152,1).Style = "Normal"
$ws1.Cells.Item(152,2).Style = "Normal"
$ws1.Cells.Item(152,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(152,2).Value = 45418.52041275083
$ws1.Cells.Item(152,3).Value = '191bugfix'
$ws1.Cells.Item(152,3).Style = "Normal"
$ws1.Cells.Item(152,4).Value = 269
$ws1.Cells.Item(152,4).Style = "Normal"
$ws1.Cells.Item(152,5).Value = 269
$ws1.Cells.Item(152,5).Style = "Normal"
$ws1.Cells.Item(152,6).Value = 0
$ws1.Cells.Item(152,6).Style = "Normal"
$ws1.Cells.Item(152,7).Value = 4.34
$ws1.Cells.Item(152,7).Style = "Normal"

$ws2 = $wb.Worksheets.Item("BETA")

# Row 47
$ws2.Cells.Item(47,1).NumberFormat = "@"
$ws2.Cells.Item(47,1).Value = '2024-04-01'
$ws2.Cells.Item(47,1).Style = "Normal"
$ws2.Cells.Item(47,2).Style = "Normal"
$ws2.Cells.Item(47,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(47,2).Value = 45383.53779894676
$ws2.Cells.Item(47,3).Value = '190betatest'
$ws2.Cells.Item(47,4).Value = 269
$ws2.Cells.Item(47,5).Value = 267
$ws2.Cells.Item(47,6).Value = 2
$ws2.Cells.Item(47,7).Value = 5.57

# Row 48
$ws2.Cells.Item(48,1).NumberFormat = "@"
$ws2.Cells.Item(48,1).Value = '2024-05-03'
$ws2.Cells.Item(48,1).Style = "Normal"
$ws2.Cells.Item(48,2).Style = "Normal"
$ws2.Cells.Item(48,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(48,2).Value = 45415.60846070602
$ws2.Cells.Item(48,3).Value = '191beta'
$ws2.Cells.Item(48,4).Value = 269
$ws2.Cells.Item(48,5).Value = 269
$ws2.Cells.Item(48,6).Value = 0
$ws2.Cells.Item(48,7).Value = 4.74

$ws3 = $wb.Worksheets.Item("AMS")

# Row 100
$ws3.Cells.Item(100,1).NumberFormat = "@"
$ws3.Cells.Item(100,1).Value = '2024-04-01'
$ws3.Cells.Item(100,1).Style = "Normal"
$ws3.Cells.Item(100,2).Style = "Normal"
$ws3.Cells.Item(100,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(100,2).Value = 45383.85344318287
$ws3.Cells.Item(100,3).Value = '190livee'
$ws3.Cells.Item(100,4).Value = 269
$ws3.Cells.Item(100,5).Value = 267
$ws3.Cells.Item(100,6).Value = 2
$ws3.Cells.Item(100,7).Value = 4.89

